$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "aaa"
$ws.Range("B2").Value = "Iansa"

# Remove row 3 entirely (previously "color" / garbled text)
$ws.Rows.Item(3).Delete()
